$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1004.1667
$ws.Cells.Item(17, 9).Value = 500
$ws.Cells.Item(17, 10).Value = 1364.2858
$ws.Cells.Item(17, 11).Value = 1500
$ws.Cells.Item(17, 12).Value = 4092.8574
$ws.Cells.Item(17, 13).Value = -1332
$ws.Cells.Item(17, 14).Value = -4428.857400000001
$ws.Cells.Item(64, 8).Value = 2982.16
$ws.Cells.Item(64, 10).Value = 3045.0588
$ws.Cells.Item(64, 12).Value = 3045.0588
$ws.Cells.Item(64, 14).Value = -3541.0588
$ws.Cells.Item(67, 8).Value = 2982.16
$ws.Cells.Item(67, 10).Value = 3045.0588
$ws.Cells.Item(67, 12).Value = 3045.0588
$ws.Cells.Item(67, 14).Value = -4761.0588
$ws.Cells.Item(74, 8).Value = 5744.2
$ws.Cells.Item(74, 9).Value = 5685.25
$ws.Cells.Item(74, 10).Value = 5980
$ws.Cells.Item(74, 11).Value = 5685.25
$ws.Cells.Item(74, 12).Value = 5980
$ws.Cells.Item(74, 13).Value = -4749.25
$ws.Cells.Item(74, 14).Value = -7852
$ws.Cells.Item(75, 8).Value = 14460.667
$ws.Cells.Item(75, 10).Value = 14460.667
$ws.Cells.Item(75, 12).Value = 14460.667
$ws.Cells.Item(75, 14).Value = -16332.667
$ws.Cells.Item(77, 8).Value = 5744.2
$ws.Cells.Item(77, 9).Value = 5685.25
$ws.Cells.Item(77, 10).Value = 5980
$ws.Cells.Item(77, 11).Value = 28426.25
$ws.Cells.Item(77, 12).Value = 29900
$ws.Cells.Item(77, 13).Value = -23746.25
$ws.Cells.Item(77, 14).Value = -39260
$ws.Cells.Item(78, 8).Value = 14460.667
$ws.Cells.Item(78, 10).Value = 14460.667
$ws.Cells.Item(78, 12).Value = 43382.001
$ws.Cells.Item(78, 14).Value = -52742.001
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).ClearContents()
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(123, 8).Value = 25000
$ws.Cells.Item(123, 10).Value = 25000
$ws.Cells.Item(123, 12).Value = 25000
$ws.Cells.Item(123, 14).Value = -34800
$ws.Cells.Item(129, 8).Value = 1850.1034
$ws.Cells.Item(129, 9).Value = 697.75
$ws.Cells.Item(129, 10).Value = 2289.0952
$ws.Cells.Item(129, 11).Value = 2093.25
$ws.Cells.Item(129, 12).Value = 6867.285600000001
$ws.Cells.Item(129, 13).Value = 2906.75
$ws.Cells.Item(129, 14).Value = -16867.2856
$ws.Cells.Item(138, 8).Value = 2788.3594
$ws.Cells.Item(138, 9).Value = 1678.862
$ws.Cells.Item(138, 10).Value = 3707.6572
$ws.Cells.Item(138, 11).Value = 5036.586
$ws.Cells.Item(138, 12).Value = 11122.9716
$ws.Cells.Item(138, 13).Value = 103.4139999999998
$ws.Cells.Item(138, 14).Value = -21402.9716
$ws.Cells.Item(141, 8).Value = 9846
$ws.Cells.Item(141, 9).Value = 2298.75
$ws.Cells.Item(141, 10).Value = 40035
$ws.Cells.Item(141, 11).Value = 6896.25
$ws.Cells.Item(141, 12).Value = 120105
$ws.Cells.Item(141, 13).Value = -1716.25
$ws.Cells.Item(141, 14).Value = -130465

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 20000
$ws.Cells.Item(62, 10).Value = 20000
$ws.Cells.Item(62, 12).Value = 20000
$ws.Cells.Item(62, 14).Value = -21248
$ws.Cells.Item(65, 8).Value = 20000
$ws.Cells.Item(65, 10).Value = 20000
$ws.Cells.Item(65, 12).Value = 60000
$ws.Cells.Item(65, 14).Value = -66240

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 14561.2
$ws.Cells.Item(82, 9).Value = 3658.8572
$ws.Cells.Item(82, 10).Value = 40000
$ws.Cells.Item(82, 11).Value = 3658.8572
$ws.Cells.Item(82, 12).Value = 40000
$ws.Cells.Item(82, 13).Value = -3275.8572
$ws.Cells.Item(82, 14).Value = -40766
$ws.Cells.Item(85, 8).Value = 14561.2
$ws.Cells.Item(85, 9).Value = 3658.8572
$ws.Cells.Item(85, 10).Value = 40000
$ws.Cells.Item(85, 11).Value = 3658.8572
$ws.Cells.Item(85, 12).Value = 40000
$ws.Cells.Item(85, 13).Value = -2332.8572
$ws.Cells.Item(85, 14).Value = -42652
$ws.Cells.Item(107, 8).Value = 998.8333
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 997.6667
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 12).Value = 997.6667
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(107, 14).Value = -4837.6667
$ws.Cells.Item(134, 8).Value = 3541.5312
$ws.Cells.Item(134, 9).Value = 2773.389
$ws.Cells.Item(134, 10).Value = 4529.143
$ws.Cells.Item(134, 11).Value = 8320.167000000001
$ws.Cells.Item(134, 12).Value = 13587.429
$ws.Cells.Item(134, 13).Value = -5785.167000000001
$ws.Cells.Item(134, 14).Value = -18657.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 758.2941
$ws.Cells.Item(5, 9).Value = 444.84
$ws.Cells.Item(5, 10).Value = 1629
$ws.Cells.Item(5, 11).Value = 1334.52
$ws.Cells.Item(5, 12).Value = 4887
$ws.Cells.Item(5, 13).Value = -1222.52
$ws.Cells.Item(5, 14).Value = -5111
$ws.Cells.Item(135, 8).Value = 758.2941
$ws.Cells.Item(135, 9).Value = 444.84
$ws.Cells.Item(135, 10).Value = 1629
$ws.Cells.Item(135, 11).Value = 4003.56
$ws.Cells.Item(135, 12).Value = 14661
$ws.Cells.Item(135, 13).Value = -1468.56
$ws.Cells.Item(135, 14).Value = -19731

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 11937.5
$ws.Cells.Item(63, 10).Value = 11937.5
$ws.Cells.Item(63, 12).Value = 11937.5
$ws.Cells.Item(63, 14).Value = -13309.5
$ws.Cells.Item(66, 8).Value = 11937.5
$ws.Cells.Item(66, 10).Value = 11937.5
$ws.Cells.Item(66, 12).Value = 35812.5
$ws.Cells.Item(66, 14).Value = -42676.5
$ws.Cells.Item(132, 8).Value = 5393.077
$ws.Cells.Item(132, 9).Value = 4181
$ws.Cells.Item(132, 10).Value = 5931.778
$ws.Cells.Item(132, 11).Value = 12543
$ws.Cells.Item(132, 12).Value = 17795.334
$ws.Cells.Item(132, 13).Value = -10013
$ws.Cells.Item(132, 14).Value = -22855.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 470.83334
$ws.Cells.Item(22, 9).Value = 468.18182
$ws.Cells.Item(22, 10).Value = 500
$ws.Cells.Item(22, 11).Value = 468.18182
$ws.Cells.Item(22, 12).Value = 500
$ws.Cells.Item(22, 13).Value = -173.18182
$ws.Cells.Item(22, 14).Value = -1090
$ws.Cells.Item(27, 8).Value = 470.83334
$ws.Cells.Item(27, 9).Value = 468.18182
$ws.Cells.Item(27, 10).Value = 500
$ws.Cells.Item(27, 11).Value = 468.18182
$ws.Cells.Item(27, 12).Value = 500
$ws.Cells.Item(27, 13).Value = -361.18182
$ws.Cells.Item(27, 14).Value = -714
$ws.Cells.Item(46, 8).Value = 357720.44
$ws.Cells.Item(46, 9).Value = 485.1875
$ws.Cells.Item(46, 10).Value = 834034.0600000001
$ws.Cells.Item(46, 11).Value = 485.1875
$ws.Cells.Item(46, 12).Value = 834034.0600000001
$ws.Cells.Item(46, 13).Value = -297.1875
$ws.Cells.Item(46, 14).Value = -834410.0600000001
$ws.Cells.Item(62, 8).Value = 5000
$ws.Cells.Item(62, 9).Value = 5000
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 5000
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -4376
$ws.Cells.Item(62, 14).Value = -6248
$ws.Cells.Item(65, 8).Value = 5000
$ws.Cells.Item(65, 9).Value = 5000
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 15000
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -11880
$ws.Cells.Item(65, 14).Value = -21240
$ws.Cells.Item(76, 8).Value = 10198.25
$ws.Cells.Item(76, 9).Value = 6000
$ws.Cells.Item(76, 10).Value = 10798
$ws.Cells.Item(76, 11).Value = 6000
$ws.Cells.Item(76, 12).Value = 10798
$ws.Cells.Item(76, 13).Value = -5662
$ws.Cells.Item(76, 14).Value = -11474
$ws.Cells.Item(79, 8).Value = 10198.25
$ws.Cells.Item(79, 9).Value = 6000
$ws.Cells.Item(79, 10).Value = 10798
$ws.Cells.Item(79, 11).Value = 6000
$ws.Cells.Item(79, 12).Value = 10798
$ws.Cells.Item(79, 13).Value = -4830
$ws.Cells.Item(79, 14).Value = -13138

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 30000
$ws.Cells.Item(64, 10).Value = 30000
$ws.Cells.Item(64, 12).Value = 30000
$ws.Cells.Item(64, 14).Value = -30496
$ws.Cells.Item(67, 8).Value = 30000
$ws.Cells.Item(67, 10).Value = 30000
$ws.Cells.Item(67, 12).Value = 30000
$ws.Cells.Item(67, 14).Value = -31716
$ws.Cells.Item(100, 8).Value = 706.73914
$ws.Cells.Item(100, 9).Value = 383.25
$ws.Cells.Item(100, 10).Value = 1446.1428
$ws.Cells.Item(100, 11).Value = 766.5
$ws.Cells.Item(100, 12).Value = 2892.2856
$ws.Cells.Item(100, 13).Value = -225.5
$ws.Cells.Item(100, 14).Value = -3974.2856
